$wb = $excel.ActiveWorkbook

# This script applies a batch of market-price/profit recalculation updates
# to the "Moogle_Profits" data sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each sheet uses columns H-N for price/profit figures that were refreshed by
# the scheduled data-update runner. Cells that should not exist after the
# update are cleared (ClearContents), and previously-empty cells that now have
# data are written directly via .Value assignment (both are value-only edits;
# no rows/columns are inserted or removed).

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 894385.5
$ws.Range("I15").Value = 894385.5
$ws.Range("K15").Value = 2683156.5
$ws.Range("M15").Value = -2682987.5
$ws.Range("H40").Value = 9958.5
$ws.Range("J40").Value = 11037.5
$ws.Range("L40").Value = 11037.5
$ws.Range("N40").Value = -11387.5
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 15000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -21240
$ws.Range("H70").Value = 2292
$ws.Range("I70").Value = 2113
$ws.Range("K70").Value = 6339
$ws.Range("M70").Value = -6069
$ws.Range("H73").Value = 2292
$ws.Range("I73").Value = 2113
$ws.Range("K73").Value = 6339
$ws.Range("M73").Value = -5403
$ws.Range("H82").Value = 9982.875
$ws.Range("I82").Value = 9982.875
$ws.Range("K82").Value = 29948.625
$ws.Range("M82").Value = -29542.625
$ws.Range("H85").Value = 9982.875
$ws.Range("I85").Value = 9982.875
$ws.Range("K85").Value = 29948.625
$ws.Range("M85").Value = -28544.625
$ws.Range("H107").Value = 725.093
$ws.Range("I107").Value = 698.6667
$ws.Range("K107").Value = 698.6667
$ws.Range("M107").Value = 1221.3333
$ws.Range("H135").Value = 1786.5667
$ws.Range("I135").Value = 1236.0454
$ws.Range("J135").Value = 3300.5
$ws.Range("K135").Value = 11124.4086
$ws.Range("L135").Value = 29704.5
$ws.Range("M135").Value = -8589.408599999999
$ws.Range("N135").Value = -34774.5
$ws.Range("H138").Value = 4353.222
$ws.Range("I138").Value = 2542.3333
$ws.Range("J138").Value = 6164.1113
$ws.Range("K138").Value = 7626.999899999999
$ws.Range("L138").Value = 18492.3339
$ws.Range("M138").Value = -2486.999899999999
$ws.Range("N138").Value = -28772.3339

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2924.7432
$ws.Range("I32").Value = 2377.8857
$ws.Range("K32").Value = 2377.8857
$ws.Range("M32").Value = -2090.8857

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 3
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 110
$ws.Range("H86").Value = 11381.2
$ws.Range("I86").Value = 3250
$ws.Range("K86").Value = 3250
$ws.Range("M86").Value = -2127
$ws.Range("H89").Value = 11381.2
$ws.Range("I89").Value = 3250
$ws.Range("K89").Value = 16250
$ws.Range("M89").Value = -10634
$ws.Range("H94").Value = 701.38464
$ws.Range("I94").Value = 510.86365
$ws.Range("K94").Value = 510.86365
$ws.Range("M94").Value = -59.86365000000001
$ws.Range("H105").Value = 5600.9653
$ws.Range("I105").Value = 4125.143
$ws.Range("J105").Value = 9475
$ws.Range("K105").Value = 4125.143
$ws.Range("L105").Value = 9475
$ws.Range("M105").Value = -2378.143
$ws.Range("N105").Value = -12969

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2433.3333
$ws.Range("I16").Value = 2500
$ws.Range("J16").Value = 2300
$ws.Range("K16").Value = 2500
$ws.Range("L16").Value = 2300
$ws.Range("M16").Value = -2213
$ws.Range("N16").Value = -2874
$ws.Range("H31").Value = 6685.8335
$ws.Range("I31").Value = 3486.1052
$ws.Range("J31").Value = 10262
$ws.Range("K31").Value = 3486.1052
$ws.Range("L31").Value = 10262
$ws.Range("M31").Value = -3191.1052
$ws.Range("N31").Value = -10852
$ws.Range("H34").Value = 6685.8335
$ws.Range("I34").Value = 3486.1052
$ws.Range("J34").Value = 10262
$ws.Range("K34").Value = 3486.1052
$ws.Range("L34").Value = 10262
$ws.Range("M34").Value = -3284.1052
$ws.Range("N34").Value = -10666
$ws.Range("H86").Value = 6280.8
$ws.Range("I86").Value = 5258.4287
$ws.Range("K86").Value = 5258.4287
$ws.Range("M86").Value = -4135.4287
$ws.Range("H89").Value = 6280.8
$ws.Range("I89").Value = 5258.4287
$ws.Range("K89").Value = 26292.1435
$ws.Range("M89").Value = -20676.1435
$ws.Range("H113").Value = 2433.3333
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = -330
$ws.Range("N113").Value = -6640
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 4450.8965
$ws.Range("I11").Value = 5574.826
$ws.Range("K11").Value = 16724.478
$ws.Range("M11").Value = -16584.478
$ws.Range("H68").Value = 1627.3422
$ws.Range("J68").Value = 1477.6765
$ws.Range("L68").Value = 4433.029500000001
$ws.Range("N68").Value = -6055.029500000001
$ws.Range("H71").Value = 1627.3422
$ws.Range("J71").Value = 1477.6765
$ws.Range("L71").Value = 13299.0885
$ws.Range("N71").Value = -21411.0885
$ws.Range("H86").Value = 1423.25
$ws.Range("I86").Value = 1712.6666
$ws.Range("K86").Value = 5137.9998
$ws.Range("M86").Value = -3951.9998
$ws.Range("H89").Value = 1423.25
$ws.Range("I89").Value = 1712.6666
$ws.Range("K89").Value = 15413.9994
$ws.Range("M89").Value = -9485.999400000001
$ws.Range("H98").Value = 2492.2
$ws.Range("I98").Value = 2006.1666
$ws.Range("J98").Value = 3221.25
$ws.Range("K98").Value = 6018.4998
$ws.Range("L98").Value = 9663.75
$ws.Range("M98").Value = -4520.4998
$ws.Range("N98").Value = -12659.75
$ws.Range("H116").Value = 2792.7
$ws.Range("I116").Value = 2255
$ws.Range("K116").Value = 6765
$ws.Range("M116").Value = -3323

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 40000
$ws.Range("I43").Value = 40000
$ws.Range("K43").Value = 40000
$ws.Range("M43").Value = -39849
$ws.Range("H63").Value = 192750
$ws.Range("I63").Value = 150000
$ws.Range("J63").Value = 199875
$ws.Range("K63").Value = 150000
$ws.Range("L63").Value = 199875
$ws.Range("M63").Value = -149314
$ws.Range("N63").Value = -201247
$ws.Range("H66").Value = 192750
$ws.Range("I66").Value = 150000
$ws.Range("J66").Value = 199875
$ws.Range("K66").Value = 450000
$ws.Range("L66").Value = 599625
$ws.Range("M66").Value = -446568
$ws.Range("N66").Value = -606489

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3926.7856
$ws.Range("J46").Value = 5999.7144
$ws.Range("L46").Value = 5999.7144
$ws.Range("N46").Value = -6375.7144
$ws.Range("H82").Value = 5323.3335
$ws.Range("I82").Value = 7241.75
$ws.Range("J82").Value = 1486.5
$ws.Range("K82").Value = 7241.75
$ws.Range("L82").Value = 1486.5
$ws.Range("M82").Value = -6880.75
$ws.Range("N82").Value = -2208.5
$ws.Range("H85").Value = 5323.3335
$ws.Range("I85").Value = 7241.75
$ws.Range("J85").Value = 1486.5
$ws.Range("K85").Value = 7241.75
$ws.Range("L85").Value = 1486.5
$ws.Range("M85").Value = -5993.75
$ws.Range("N85").Value = -3982.5
$ws.Range("H136").Value = 7139.6606
$ws.Range("I136").Value = 4241.05
$ws.Range("J136").Value = 8750
$ws.Range("K136").Value = 12723.15
$ws.Range("L136").Value = 26250
$ws.Range("M136").Value = -10173.15
$ws.Range("N136").Value = -31350

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 413.84616
$ws.Range("J113").Value = 503.8
$ws.Range("L113").Value = 1511.4
$ws.Range("N113").Value = -5851.4
$ws.Range("H136").Value = 3746.7754
$ws.Range("I136").Value = 3288.762
$ws.Range("K136").Value = 9866.286
$ws.Range("M136").Value = -7316.286
